$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data range (B2:C5) since the table is being rebuilt in A1:C7
$ws.Range("B2:C5").ClearContents()

# Header row
$ws.Range("A1").Value = "建物"
$ws.Range("B1").Value = "面積"
$ws.Range("C1").Value = "竣工年"

# Data rows
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 2020

$ws.Range("A3").Value = "B"
$ws.Range("B3").Value = 2500
$ws.Range("C3").Value = 2015

$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = 1800
$ws.Range("C4").Value = 2018

$ws.Range("A5").Value = "D"
$ws.Range("B5").Value = 5600
$ws.Range("C5").Value = 2000

$ws.Range("A6").Value = "E"
$ws.Range("B6").Value = 4700
$ws.Range("C6").Value = 1995

$ws.Range("A7").Value = "F"
$ws.Range("B7").Value = 8100
$ws.Range("C7").Value = 2013

# Update selection to mirror the target file (active cell C8)
$ws.Range("C8").Select()
